$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold font, border, date number format) from A92 to the new date cells
$ws.Range("A92").Copy($ws.Range("A93:A98"))

# Row 93
$ws.Cells.Item(93, 1).Value = 45523
$ws.Range("B93").Value = 755.707627235
$ws.Range("C93").Value = 182.795579686
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("I93").Value = 258.343034899
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0.00958560036
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("N93").Value = 106.0313597808
$ws.Range("O93").Value = 56.803971504
$ws.Range("P93").Value = 0
$ws.Range("Q93").Value = 0.0000018432
$ws.Range("R93").Value = 0
$ws.Range("S93").Value = 0
$ws.Range("T93").Value = 0
$ws.Range("U93").Value = 268.1392220048816
$ws.Range("W93").Value = 0
$ws.Range("X93").Value = 0
$ws.Range("Y93").Value = 0
$ws.Range("Z93").Value = 183.865549291598

# Row 94
$ws.Cells.Item(94, 1).Value = 45524
$ws.Range("B94").Value = 750.307944718
$ws.Range("C94").Value = 178.389948007
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 0
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("I94").Value = 254.726339837
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0.01233964356
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 0
$ws.Range("N94").Value = 102.81524270688
$ws.Range("O94").Value = 57.922960914
$ws.Range("P94").Value = 0
$ws.Range("Q94").Value = 0.00000186
$ws.Range("R94").Value = 0
$ws.Range("S94").Value = 0
$ws.Range("T94").Value = 0
$ws.Range("U94").Value = 262.1265581526729
$ws.Range("W94").Value = 0
$ws.Range("X94").Value = 0
$ws.Range("Y94").Value = 0
$ws.Range("Z94").Value = 176.666558677714

# Row 95
$ws.Cells.Item(95, 1).Value = 45525
$ws.Range("B95").Value = 777.5444925833
$ws.Range("C95").Value = 182.4038293085
$ws.Range("D95").Value = 0
$ws.Range("E95").Value = 0
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("I95").Value = 255.657370249
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0.013245226515
$ws.Range("L95").Value = 0
$ws.Range("M95").Value = 0
$ws.Range("N95").Value = 104.28929636576
$ws.Range("O95").Value = 57.994169331
$ws.Range("P95").Value = 0
$ws.Range("Q95").Value = 0.0000018984
$ws.Range("R95").Value = 0
$ws.Range("S95").Value = 0
$ws.Range("T95").Value = 0
$ws.Range("U95").Value = 265.1968545878432
$ws.Range("W95").Value = 0
$ws.Range("X95").Value = 0
$ws.Range("Y95").Value = 0
$ws.Range("Z95").Value = 186.570372354618

# Row 96
$ws.Cells.Item(96, 1).Value = 45526
$ws.Range("B96").Value = 767.6250711023999
$ws.Range("C96").Value = 181.860925688
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("I96").Value = 256.427260782
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0.0148603581
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("N96").Value = 107.23740368352
$ws.Range("O96").Value = 59.296266099
$ws.Range("P96").Value = 0
$ws.Range("Q96").Value = 0.0000019176
$ws.Range("R96").Value = 0
$ws.Range("S96").Value = 0
$ws.Range("T96").Value = 0
$ws.Range("U96").Value = 279.3969756005064
$ws.Range("W96").Value = 0
$ws.Range("X96").Value = 0
$ws.Range("Y96").Value = 0
$ws.Range("Z96").Value = 187.797945898604

# Row 97
$ws.Cells.Item(97, 1).Value = 45527
$ws.Range("B97").Value = 814.1765134564
$ws.Range("C97").Value = 191.540280148
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("I97").Value = 274.403309506
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0.016651442445
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("N97").Value = 120.2693780768
$ws.Range("O97").Value = 60.19145762700001
$ws.Range("P97").Value = 0
$ws.Range("Q97").Value = 0.0000021792
$ws.Range("R97").Value = 0
$ws.Range("S97").Value = 0
$ws.Range("T97").Value = 0
$ws.Range("U97").Value = 296.9232510846041
$ws.Range("W97").Value = 0
$ws.Range("X97").Value = 0
$ws.Range("Y97").Value = 0
$ws.Range("Z97").Value = 209.145241765208

# Row 98
$ws.Cells.Item(98, 1).Value = 45528
$ws.Range("B98").Value = 815.6992824111001
$ws.Range("C98").Value = 191.9230168
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("I98").Value = 287.670492877
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 58.91118216102
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("N98").Value = 120.83889880864
$ws.Range("O98").Value = 0.128122155
$ws.Range("P98").Value = 0
$ws.Range("Q98").Value = 0.00000228
$ws.Range("R98").Value = 0
$ws.Range("S98").Value = 0
$ws.Range("T98").Value = 0
$ws.Range("W98").Value = 0
$ws.Range("X98").Value = 0
$ws.Range("Y98").Value = 0
$ws.Range("Z98").Value = 212.078934472022
